$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3525
$ws.Range("I76").Value = 3542.8572
$ws.Range("J76").Value = 3400
$ws.Range("K76").Value = 3542.8572
$ws.Range("L76").Value = 3400
$ws.Range("M76").Value = -3227.8572
$ws.Range("N76").Value = -4030
$ws.Range("H79").Value = 3525
$ws.Range("I79").Value = 3542.8572
$ws.Range("J79").Value = 3400
$ws.Range("K79").Value = 3542.8572
$ws.Range("L79").Value = 3400
$ws.Range("M79").Value = -2450.8572
$ws.Range("N79").Value = -5584
$ws.Range("H135").Value = 37037.25
$ws.Range("I135").Value = 42805.293
$ws.Range("J135").Value = 2429
$ws.Range("K135").Value = 385247.637
$ws.Range("L135").Value = 21861
$ws.Range("M135").Value = -382712.637
$ws.Range("N135").Value = -26931
$ws.Range("H137").Value = 3573397
$ws.Range("I137").Value = 4349582
$ws.Range("K137").Value = 13048746
$ws.Range("M137").Value = -13046196
$ws.Range("H138").Value = 2876353.5
$ws.Range("I138").Value = 1038.2222
$ws.Range("J138").Value = 4170245.2
$ws.Range("K138").Value = 3114.6666
$ws.Range("L138").Value = 12510735.6
$ws.Range("M138").Value = 2025.3334
$ws.Range("N138").Value = -12521015.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 33358.758
$ws.Range("I132").Value = 22960.756
$ws.Range("J132").Value = 60882.883
$ws.Range("K132").Value = 68882.26800000001
$ws.Range("L132").Value = 182648.649
$ws.Range("M132").Value = -66352.26800000001
$ws.Range("N132").Value = -187708.649
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1804.75
$ws.Range("I107").Value = 1841.5454
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 1841.5454
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = 78.45460000000003
$ws.Range("N107").Value = -5240
$ws.Range("H134").Value = 2329.4407
$ws.Range("I134").Value = 1814.7826
$ws.Range("J134").Value = 4150.5386
$ws.Range("K134").Value = 5444.3478
$ws.Range("L134").Value = 12451.6158
$ws.Range("M134").Value = -2909.3478
$ws.Range("N134").Value = -17521.6158
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1240.0714
$ws.Range("I16").Value = 865.8570999999999
$ws.Range("J16").Value = 1614.2858
$ws.Range("K16").Value = 865.8570999999999
$ws.Range("L16").Value = 1614.2858
$ws.Range("M16").Value = -578.8570999999999
$ws.Range("N16").Value = -2188.2858
$ws.Range("H31").Value = 2493.7908
$ws.Range("I31").Value = 1359.4839
$ws.Range("J31").Value = 5424.0835
$ws.Range("K31").Value = 1359.4839
$ws.Range("L31").Value = 5424.0835
$ws.Range("M31").Value = -1064.4839
$ws.Range("N31").Value = -6014.0835
$ws.Range("H34").Value = 2493.7908
$ws.Range("I34").Value = 1359.4839
$ws.Range("J34").Value = 5424.0835
$ws.Range("K34").Value = 1359.4839
$ws.Range("L34").Value = 5424.0835
$ws.Range("M34").Value = -1157.4839
$ws.Range("N34").Value = -5828.0835
$ws.Range("H105").Value = 1108.6428
$ws.Range("I105").Value = 1101.1111
$ws.Range("J105").Value = 1122.2
$ws.Range("K105").Value = 1101.1111
$ws.Range("L105").Value = 1122.2
$ws.Range("M105").Value = 645.8888999999999
$ws.Range("N105").Value = -4616.2
$ws.Range("H113").Value = 1240.0714
$ws.Range("I113").Value = 865.8570999999999
$ws.Range("J113").Value = 1614.2858
$ws.Range("K113").Value = 865.8570999999999
$ws.Range("L113").Value = 1614.2858
$ws.Range("M113").Value = 1304.1429
$ws.Range("N113").Value = -5954.2858
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 599.7954999999999
$ws.Range("I113").Value = 495.83334
$ws.Range("J113").Value = 638.78125
$ws.Range("K113").Value = 1487.50002
$ws.Range("L113").Value = 1916.34375
$ws.Range("M113").Value = 682.4999800000001
$ws.Range("N113").Value = -6256.34375
$ws.Range("H131").Value = 978.625
$ws.Range("I131").Value = 425.8
$ws.Range("J131").Value = 1081
$ws.Range("K131").Value = 1277.4
$ws.Range("L131").Value = 3243
$ws.Range("M131").Value = 3762.6
$ws.Range("N131").Value = -13323
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1593.3334
$ws.Range("I126").Value = 1418.1818
$ws.Range("J126").Value = 2075
$ws.Range("K126").Value = 4254.5454
$ws.Range("L126").Value = 6225
$ws.Range("M126").Value = -1784.5454
$ws.Range("N126").Value = -11165
$ws.Range("H132").Value = 38406.203
$ws.Range("I132").Value = 26160
$ws.Range("J132").Value = 73395.36
$ws.Range("K132").Value = 78480
$ws.Range("L132").Value = 220186.08
$ws.Range("M132").Value = -75950
$ws.Range("N132").Value = -225246.08
$ws.Range("H135").Value = 31197.777
$ws.Range("J135").Value = 31197.777
$ws.Range("L135").Value = 31197.777
$ws.Range("N135").Value = -41337.777
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 96536
$ws.Range("J139").Value = 96536
$ws.Range("L139").Value = 96536
$ws.Range("N139").Value = -106816
$ws.Range("H140").Value = 500000
$ws.Range("J140").Value = 500000
$ws.Range("L140").Value = 500000
$ws.Range("N140").Value = -510360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2769.1667
$ws.Range("I61").Value = 2769.1667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2769.1667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2567.1667
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2769.1667
$ws.Range("I113").Value = 2769.1667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2769.1667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -599.1667000000002
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 41000.28
$ws.Range("I132").Value = 17937.906
$ws.Range("J132").Value = 146428.28
$ws.Range("K132").Value = 53813.71799999999
$ws.Range("L132").Value = 439284.84
$ws.Range("M132").Value = -51283.71799999999
$ws.Range("N132").Value = -444344.84
$ws.Range("H136").Value = 61483.65
$ws.Range("I136").Value = 32513.875
$ws.Range("J136").Value = 525000
$ws.Range("K136").Value = 97541.625
$ws.Range("L136").Value = 1575000
$ws.Range("M136").Value = -94991.625
$ws.Range("N136").Value = -1580100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 467.875
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 523.8333
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1571.4999
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -5411.4999
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 55526.27
$ws.Range("I132").Value = 40987.16
$ws.Range("J132").Value = 85816.086
$ws.Range("K132").Value = 122961.48
$ws.Range("L132").Value = 257448.258
$ws.Range("M132").Value = -120431.48
$ws.Range("N132").Value = -262508.258
